$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.78%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.48%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.064"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.85%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07660"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.56%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.604"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.64%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9971"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "8.11%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1240"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.88%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1865"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.63%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09044"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.48%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04170"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.97%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1046"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.71%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001280"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.03%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.005765"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.88%"
$ws.Range("B15").Value = "UpBots"
$ws.Range("C15").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.007430"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1,891.71%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.327"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.25%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.418"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.36%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.356"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.89%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3351"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.76%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.432"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "6.07%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1400"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.21%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3198"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "13.83%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04152"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.16%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001282"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.88%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004491"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "18.44%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.36%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02444"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.27%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05277"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.86%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005969"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.80%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007664"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.25%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1345"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.00%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007344"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.29%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007453"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.68%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3014"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.14%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006722"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "8.61%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.25%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04436"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-4.24%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.22%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.25%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.25%"
